# Journal TP3 - "Premiere entree du journal" commit.
# The document starts as a single empty paragraph. We rebuild a 13
# paragraph skeleton (keeping that original empty paragraph in place as
# paragraph #2) and then fill each paragraph with its text/formatting.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Build the paragraph skeleton.
# ---------------------------------------------------------------------
$orig = $d.Paragraphs.Item(1)
$orig.Range.InsertParagraphBefore()
# Paragraphs.Item(1) = brand new empty paragraph (will hold the title)
# Paragraphs.Item(2) = the original, untouched empty paragraph

for ($i = 0; $i -lt 11; $i++) {
    $last = $d.Paragraphs.Item($d.Paragraphs.Count)
    $last.Range.InsertParagraphAfter()
}
# Now there are 13 paragraphs total:
#  1  Title
#  2  (blank - original)
#  3  Nom / Da line
#  4  (blank)
#  5  (blank)
#  6  Premiere entree :
#  7  (blank)
#  8  Date de la journee : ...
#  9  Probleme rencontre : ...
# 10  (blank, carries formatting)
# 11  Solution trouvee : ...
# 12  (blank, carries formatting)
# 13  Apprentissage : ...

# ---------------------------------------------------------------------
# Helper text.
# ---------------------------------------------------------------------
$TITLE               = "Journal TP3-Applications-Web"
$NOM                 = "Nom : Laurier Mainguy                         Da : 6228273"
$PREMIERE_ENTREE     = "Première entrée : "
$DATE_LABEL          = "Date de la journée"
$DATE_REST           = " : 22 mai 2024"
$PROBLEME_LABEL      = "Problème rencontré"
$PROBLEME_REST       = " : Je me demandais comment j’allais introduire une procédure de livraison avec seulement les notions qu’on a appris en cours, car nous n’avons pas appris à modifier des données du site à la base de données."
$SOLUTION_LABEL      = "Solution trouvée"
$SOLUTION_REST       = " : J’ai décidé que j’allais juste afficher un message qui remercie au client d’avoir commandé un tel produit. Je vais me servir de fonction javascript pour faire en sorte de réafficher les données du formulaire que l’utilisateur à entrer. Grâce à cela notre site web sera dynamique."
$APPRENTISSAGE_LABEL = "Apprentissage :"
$APPRENTISSAGE_REST  = " Comment utiliser Javascript pour l’affichage de données dynamiques."

# ---------------------------------------------------------------------
# 2) Paragraph 1 - Title.
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(1)
$p.Range.Text = $TITLE
$r = $p.Range
$r.Font.Name = "Times New Roman"
$r.Font.NameAscii = "Times New Roman"
$r.Font.NameBi = "Times New Roman"
$r.Font.Bold = $true
$r.Font.BoldBi = $true
$r.Font.Size = 28
$r.Font.SizeBi = 28

# ---------------------------------------------------------------------
# 3) Paragraph 3 - Nom / Da.
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(3)
$p.Range.Text = $NOM
$r = $p.Range
$r.Font.Name = "Times New Roman"
$r.Font.NameAscii = "Times New Roman"
$r.Font.NameBi = "Times New Roman"
$r.Font.Bold = $true
$r.Font.BoldBi = $true
$r.Font.Size = 14
$r.Font.SizeBi = 14

# ---------------------------------------------------------------------
# 4) Paragraph 6 - "Premiere entree : " (bold, default font/size).
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(6)
$p.Range.Text = $PREMIERE_ENTREE
$r = $p.Range
$r.Font.Bold = $true
$r.Font.BoldBi = $true

# ---------------------------------------------------------------------
# 5) Paragraph 8 - "Date de la journee" (bold) + rest (regular).
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(8)
$fullText = $DATE_LABEL + $DATE_REST
$p.Range.Text = $fullText
$pr = $p.Range
$pr.Font.Name = "Times New Roman"
$pr.Font.NameAscii = "Times New Roman"
$pr.Font.NameBi = "Times New Roman"
$pr.Font.Size = 12
$pr.Font.SizeBi = 12
$labelStart = $p.Range.Start
$labelEnd = $labelStart + $DATE_LABEL.Length
$labelRange = $d.Range($labelStart, $labelEnd)
$labelRange.Font.Bold = $true

# ---------------------------------------------------------------------
# 6) Paragraph 9 - "Probleme rencontre" (bold) + rest (regular).
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(9)
$fullText = $PROBLEME_LABEL + $PROBLEME_REST
$p.Range.Text = $fullText
$pr = $p.Range
$pr.Font.Name = "Times New Roman"
$pr.Font.NameAscii = "Times New Roman"
$pr.Font.NameBi = "Times New Roman"
$pr.Font.Size = 12
$pr.Font.SizeBi = 12
$labelStart = $p.Range.Start
$labelEnd = $labelStart + $PROBLEME_LABEL.Length
$labelRange = $d.Range($labelStart, $labelEnd)
$labelRange.Font.Bold = $true

# ---------------------------------------------------------------------
# 7) Paragraph 10 - blank, carries Times New Roman / 12pt formatting.
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(10)
$pr = $p.Range
$pr.Font.Name = "Times New Roman"
$pr.Font.NameAscii = "Times New Roman"
$pr.Font.NameBi = "Times New Roman"
$pr.Font.Size = 12
$pr.Font.SizeBi = 12

# ---------------------------------------------------------------------
# 8) Paragraph 11 - "Solution trouvee" (bold) + rest (regular).
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(11)
$fullText = $SOLUTION_LABEL + $SOLUTION_REST
$p.Range.Text = $fullText
$pr = $p.Range
$pr.Font.Name = "Times New Roman"
$pr.Font.NameAscii = "Times New Roman"
$pr.Font.NameBi = "Times New Roman"
$pr.Font.Size = 12
$pr.Font.SizeBi = 12
$labelStart = $p.Range.Start
$labelEnd = $labelStart + $SOLUTION_LABEL.Length
$labelRange = $d.Range($labelStart, $labelEnd)
$labelRange.Font.Bold = $true

# ---------------------------------------------------------------------
# 9) Paragraph 12 - blank, carries Times New Roman / 12pt formatting.
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(12)
$pr = $p.Range
$pr.Font.Name = "Times New Roman"
$pr.Font.NameAscii = "Times New Roman"
$pr.Font.NameBi = "Times New Roman"
$pr.Font.Size = 12
$pr.Font.SizeBi = 12

# ---------------------------------------------------------------------
# 10) Paragraph 13 - "Apprentissage :" (bold) + rest (regular).
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(13)
$fullText = $APPRENTISSAGE_LABEL + $APPRENTISSAGE_REST
$p.Range.Text = $fullText
$pr = $p.Range
$pr.Font.Name = "Times New Roman"
$pr.Font.NameAscii = "Times New Roman"
$pr.Font.NameBi = "Times New Roman"
$pr.Font.Size = 12
$pr.Font.SizeBi = 12
$labelStart = $p.Range.Start
$labelEnd = $labelStart + $APPRENTISSAGE_LABEL.Length
$labelRange = $d.Range($labelStart, $labelEnd)
$labelRange.Font.Bold = $true

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
